$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("110:111").Insert()
$ws.Range("C139").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C141"), "https://github.com/ilanc/flot/")
